# The workbook's stock-report rows were generated in pairs (two lots of the
# same item), and each pair had its B:G (Code, Name, Rate, Value, Qty,
# Amount) values swapped between the two rows - i.e. the two rows in every
# pair traded places while the serial number in column A stayed put.
#
# Row pairs affected (1-based worksheet rows):
#   151/152, 198/199, 228/229, 237/238, 326/327, 371/372, 373/374,
#   387/388, 391/392, 401/402, 484/485, 554/555, 560/561, 563/564,
#   568/569, 644/645

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(151, 152),
    @(198, 199),
    @(228, 229),
    @(237, 238),
    @(326, 327),
    @(371, 372),
    @(373, 374),
    @(387, 388),
    @(391, 392),
    @(401, 402),
    @(484, 485),
    @(554, 555),
    @(560, 561),
    @(563, 564),
    @(568, 569),
    @(644, 645)
)

# Columns B through G (2..7) hold Code, Name, Rate, Value, Qty, Amount.
$cols = 2, 3, 4, 5, 6, 7

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}
